$wb = $excel.ActiveWorkbook

# --- Add the new "Scan Path Switching" worksheet after Sheet1 ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Scan Path Switching"

# Populate cells in the exact order needed to reproduce the shared-string table.
# Headers (min_x..max_z)
$ws.Range("A1").Value = "min_x"
$ws.Range("B1").Value = "min_y"
$ws.Range("C1").Value = "min_z"
$ws.Range("D1").Value = "max_x"
$ws.Range("E1").Value = "max_y"
$ws.Range("F1").Value = "max_z"

# Notes column entered next
$ws.Range("I3").Value = "Note: Capitalization Matters"
$ws.Range("I4").Value = "Note: Any unspecified areas will have ``default`` hatching applied to them"

# scanpath header + eligible-values note
$ws.Range("G1").Value = "scanpath"
$ws.Range("I2").Value = "Eligible ``scanpath`` Values: ``default``, ``island``"

# scanpath values
$ws.Range("G2").Value = "island"
$ws.Range("G3").Value = "island"
$ws.Range("G4").Value = "frick"

# Numeric grid
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 6

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 7

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 8

# Column width for G (target stored width 17.5703125; engine quantizes to 1/6 px units)
$ws.Columns.Item(7).ColumnWidth = 16.666666666666668

# Selection on new sheet
$ws.Range("G5").Select()

# --- Update Sheet1 selection ---
$sheet1.Range("D14").Select()

# --- Make the new sheet active tab ---
$ws.Activate()
